# Add season-record columns (Wins, Losses, Ties) to the stats table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header row: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties" --
$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Copy the formatting of the existing header cell (A1) onto the new
# header cells so they pick up the same bold/centered/bordered style
# used by every other header (maps to cellXfs index 1 / s="1").
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# -- Data rows 2-50: Wins = 80, Losses = 82, Ties = 0 --
for ($row = 2; $row -le 50; $row++) {
    $ws.Cells.Item($row, 30).Value = 80
    $ws.Cells.Item($row, 31).Value = 82
    $ws.Cells.Item($row, 32).Value = 0
}
